$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36/37: swap Monero <-> Fetch.AI, with updated price/volume data ---
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

# --- Price column (D): force text type (matches source inlineStr) without altering cell style ---
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D48","D49","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '65.493.22'
$ws.Range("D3").Value = '3.433.25'
$ws.Range("D5").Value = '594.78'
$ws.Range("D6").Value = '135.07'
$ws.Range("D7").Value = '3.432.93'
$ws.Range("D8").Value = '0.998'
$ws.Range("D9").Value = '0.491'
$ws.Range("D10").Value = '7.47'
$ws.Range("D11").Value = '0.122'
$ws.Range("D12").Value = '0.377'
$ws.Range("D13").Value = '4.005.25'
$ws.Range("D14").Value = '0.0000180'
$ws.Range("D15").Value = '26.34'
$ws.Range("D16").Value = '3.438.53'
$ws.Range("D17").Value = '65.342.28'
$ws.Range("D19").Value = '9.98'
$ws.Range("D20").Value = '5.72'
$ws.Range("D21").Value = '13.66'
$ws.Range("D22").Value = '391.82'
$ws.Range("D23").Value = '73.33'
$ws.Range("D24").Value = '0.543'
$ws.Range("D25").Value = '0.999'
$ws.Range("D26").Value = '3.574.22'
$ws.Range("D27").Value = '0.0000105'
$ws.Range("D28").Value = '0.998'
$ws.Range("D29").Value = '2.25'
$ws.Range("D30").Value = '7.15'
$ws.Range("D31").Value = '8.18'
$ws.Range("D32").Value = '3.438.20'
$ws.Range("D34").Value = '0.145'
$ws.Range("D35").Value = '22.57'
$ws.Range("D38").Value = '6.83'
$ws.Range("D40").Value = '4.83'
$ws.Range("D41").Value = '0.0769'
$ws.Range("D42").Value = '0.811'
$ws.Range("D43").Value = '43.55'
$ws.Range("D44").Value = '0.999'
$ws.Range("D45").Value = '4.39'
$ws.Range("D48").Value = '21.78'
$ws.Range("D49").Value = '6.54'
$ws.Range("D50").Value = '2.13'
$ws.Range("D51").Value = '2.193.95'
$ws.Range("D36").Value = '1.23'
$ws.Range("D37").Value = '172.22'

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# --- Volume(1h) column (E): plain text percentage strings ---
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("E3").Value = '  -4.51%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("E6").Value = '  -9.05%  '
$ws.Range("E7").Value = '  -4.50%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("E10").Value = '  -5.33%  '
$ws.Range("E11").Value = '  -10.26%  '
$ws.Range("E12").Value = '  -8.97%  '
$ws.Range("E13").Value = '  -4.73%  '
$ws.Range("E14").Value = '  -12.38%  '
$ws.Range("E15").Value = '  -10.90%  '
$ws.Range("E16").Value = '  -4.39%  '
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("E19").Value = '  -10.30%  '
$ws.Range("E20").Value = '  -9.56%  '
$ws.Range("E21").Value = '  -8.44%  '
$ws.Range("E22").Value = '  -7.40%  '
$ws.Range("E23").Value = '  -6.82%  '
$ws.Range("E24").Value = '  -11.25%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("E27").Value = '  -12.34%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("E29").Value = '  -9.83%  '
$ws.Range("E30").Value = '  -13.26%  '
$ws.Range("E31").Value = '  -12.83%  '
$ws.Range("E32").Value = '  -4.29%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -7.93%  '
$ws.Range("E35").Value = '  -10.13%  '
$ws.Range("E38").Value = '  -12.00%  '
$ws.Range("E39").Value = '  -9.14%  '
$ws.Range("E40").Value = '  -13.59%  '
$ws.Range("E41").Value = '  -9.66%  '
$ws.Range("E42").Value = '  -8.31%  '
$ws.Range("E43").Value = '  -5.16%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("E45").Value = '  -15.45%  '
$ws.Range("E46").Value = '  -13.07%  '
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("E48").Value = '  -8.08%  '
$ws.Range("E49").Value = '  -8.49%  '
$ws.Range("E50").Value = '  -15.74%  '
$ws.Range("E51").Value = '  -8.32%  '
$ws.Range("E36").Value = '  -14.22%  '
$ws.Range("E37").Value = '  -1.55%  '

